# Weekly data refresh: insert one new observation for "Ají" (Macroferia
# Regional de Talca) above the existing row 268, shifting all subsequent
# rows down by one. This mirrors the target diff, where a brand-new row
# of data appears at row 268 and the previous rows 268-361 become 269-362.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 268, pushing rows 268:361 down to 269:362.
$ws.Rows("268:268").Insert()

# Populate the newly inserted row 268 with the new weekly observation.
$ws.Range("A268").Value = 5
$ws.Range("B268").Value = "Macroferia Regional de Talca"
$ws.Range("C268").Value = "Maule"
$ws.Range("D268").Value = 45027
$ws.Range("E268").Value = 7
$ws.Range("F268").Value = 100112021
$ws.Range("G268").Value = "Ají"
$ws.Range("H268").Value = "Cristal"
$ws.Range("I268").Value = "Primera"
$ws.Range("J268").Value = 100
$ws.Range("K268").Value = 15000
$ws.Range("L268").Value = 15000
$ws.Range("M268").Value = 15000
$ws.Range("N268").Value = '$/saco 25 kilos'
$ws.Range("O268").Value = "Región del Maule"
$ws.Range("P268").Value = 600
$ws.Range("Q268").Value = 25
$ws.Range("R268").Value = "Hortaliza"
